$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Table 1: 항공 (flight) block — header at row 4, new "English name" cell at
# C3 and a new data-type annotation row inserted at row 5 (was blank, so a
# simple write is enough — nothing below needs to shift).
# ---------------------------------------------------------------------------
$ws.Range("C3").Value = "fligth"

$ws.Range("B5").Value = "varchar"
$ws.Range("C5").Value = "varchar"
$ws.Range("D5").Value = "date"
$ws.Range("E5").Value = "varchar"
$ws.Range("F5").Value = "varchar"
$ws.Range("G5").Value = "number"
$ws.Range("H5").Value = "number"
$ws.Range("I5").Value = "varchar"
$ws.Range("J5").Value = "varcahr"
$ws.Range("K5").Value = "varchar"

# ---------------------------------------------------------------------------
# Table 2: 승객 (passenger) block — header at row 14, new "English name"
# cell at C13 and a new data-type row inserted at row 15 (also blank
# beforehand).
# ---------------------------------------------------------------------------
$ws.Range("C13").Value = "passenger"

$ws.Range("B15").Value = "varchar"
$ws.Range("C15").Value = "varchar"
$ws.Range("D15").Value = "number"
$ws.Range("F15").Value = "varchar"
$ws.Range("G15").Value = "varchar"
$ws.Range("H15").Value = "number"

# ---------------------------------------------------------------------------
# Table 3: 나라 (nation) block — header at row 23, data already present in
# rows 24-25, so inserting the new data-type row has to shift the existing
# data down by one first (row25->26, row24->25) before writing the new
# row24.
# ---------------------------------------------------------------------------
$ws.Range("C22").Value = "nation"

$ws.Range("B26").Value = $ws.Range("B25").Value2
$ws.Range("C26").Value = $ws.Range("C25").Value2
$ws.Range("B25").Value = $ws.Range("B24").Value2
$ws.Range("C25").Value = $ws.Range("C24").Value2

$ws.Range("B24").Value = "varchar"
$ws.Range("C24").Value = "varchar"

# ---------------------------------------------------------------------------
# Table 4: 공항 (airport) block — header at row 31, data already present in
# rows 32-33, so the same shift-then-insert pattern applies (row33->34,
# row32->33) before writing the new row32.
# ---------------------------------------------------------------------------
$ws.Range("C30").Value = "airport"

$ws.Range("B34").Value = $ws.Range("B33").Value2
$ws.Range("C34").Value = $ws.Range("C33").Value2
$ws.Range("D34").Value = $ws.Range("D33").Value2
$ws.Range("E34").Value = $ws.Range("E33").Value2

$ws.Range("B33").Value = $ws.Range("B32").Value2
$ws.Range("C33").Value = $ws.Range("C32").Value2
$ws.Range("D33").Value = $ws.Range("D32").Value2
$ws.Range("E33").Value = $ws.Range("E32").Value2

$ws.Range("B32").Value = "varchar"
$ws.Range("C32").Value = "varchar"
$ws.Range("D32").Value = "varchar"
$ws.Range("E32").Value = "varchar"

# ---------------------------------------------------------------------------
# Selection moves to A6 (as captured in the saved view state).
# ---------------------------------------------------------------------------
$ws.Range("A6").Select()
